$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Range, [string]$Text)
    $Range.Value = "'" + $Text
    $Range.Style = "Normal"
}

Set-TextCell $ws.Range("D2") '43.111.42'
Set-TextCell $ws.Range("E2") '  +0.30%  '
Set-TextCell $ws.Range("E3") '  +0.28%  '
Set-TextCell $ws.Range("E4") '  +0.09%  '
Set-TextCell $ws.Range("D5") '302.15'
Set-TextCell $ws.Range("E5") '  -0.20%  '
Set-TextCell $ws.Range("D6") '98.93'
Set-TextCell $ws.Range("E6") '  -1.98%  '
Set-TextCell $ws.Range("D7") '0.511'
Set-TextCell $ws.Range("E7") '  +1.13%  '
Set-TextCell $ws.Range("E8") '  +0.01%  '
Set-TextCell $ws.Range("D9") '0.523'
Set-TextCell $ws.Range("E9") '  +0.82%  '
Set-TextCell $ws.Range("D10") '35.82'
Set-TextCell $ws.Range("E10") '  +1.48%  '
Set-TextCell $ws.Range("D11") '0.0791'
Set-TextCell $ws.Range("E11") '  -0.65%  '
Set-TextCell $ws.Range("E12") '  -1.03%  '
Set-TextCell $ws.Range("D13") '17.90'
Set-TextCell $ws.Range("E13") '  -0.29%  '
Set-TextCell $ws.Range("D14") '6.92'
Set-TextCell $ws.Range("E14") '  +0.03%  '
Set-TextCell $ws.Range("D15") '2.674.08'
Set-TextCell $ws.Range("D16") '2.292.24'
Set-TextCell $ws.Range("E16") '  +0.97%  '
Set-TextCell $ws.Range("E17") '  -2.71%  '
Set-TextCell $ws.Range("D18") '43.014.93'
Set-TextCell $ws.Range("E18") '  +0.29%  '
Set-TextCell $ws.Range("D19") '13.57'
Set-TextCell $ws.Range("E19") '  +7.32%  '
Set-TextCell $ws.Range("D20") '0.0₃0912'
Set-TextCell $ws.Range("E20") '  +0.68%  '
Set-TextCell $ws.Range("D21") '6.19'
Set-TextCell $ws.Range("E21") '  +0.16%  '
Set-TextCell $ws.Range("D22") '68.10'
Set-TextCell $ws.Range("E22") '  +0.22%  '
Set-TextCell $ws.Range("D23") '240.61'
Set-TextCell $ws.Range("E23") '  +1.41%  '
Set-TextCell $ws.Range("D24") '2.18'
Set-TextCell $ws.Range("E24") '  -0.68%  '
Set-TextCell $ws.Range("E25") '  -0.06%  '
Set-TextCell $ws.Range("D26") '2.45'
Set-TextCell $ws.Range("E26") '  -0.84%  '
Set-TextCell $ws.Range("D27") '24.98'
Set-TextCell $ws.Range("E27") '  +0.78%  '
Set-TextCell $ws.Range("D28") '168.42'
Set-TextCell $ws.Range("E28") '  +0.61%  '
Set-TextCell $ws.Range("D29") '9.19'
Set-TextCell $ws.Range("E29") '  -0.55%  '
Set-TextCell $ws.Range("D30") '2.04'
Set-TextCell $ws.Range("E30") '  -1.63%  '
Set-TextCell $ws.Range("D31") '33.50'
Set-TextCell $ws.Range("E31") '  -1.78%  '
Set-TextCell $ws.Range("B32") 'RenderToken'
Set-TextCell $ws.Range("C32") 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell $ws.Range("D32") '4.97'
Set-TextCell $ws.Range("E32") '  +7.32%  '
Set-TextCell $ws.Range("B33") 'Filecoin'
Set-TextCell $ws.Range("C33") 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws.Range("D33") '5.24'
Set-TextCell $ws.Range("E33") '  +4.24%  '
Set-TextCell $ws.Range("B34") 'Celestia'
Set-TextCell $ws.Range("C34") 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextCell $ws.Range("D34") '18.44'
Set-TextCell $ws.Range("E34") '  +8.57%  '
Set-TextCell $ws.Range("B35") 'FirstDigitalUSD'
Set-TextCell $ws.Range("C35") 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextCell $ws.Range("D35") '0.999'
Set-TextCell $ws.Range("E35") '  +0.00%  '
Set-TextCell $ws.Range("E36") '  -0.12%  '
Set-TextCell $ws.Range("E37") '  +0.49%  '
Set-TextCell $ws.Range("E38") '  +0.34%  '
Set-TextCell $ws.Range("E39") '  +0.94%  '
Set-TextCell $ws.Range("D40") '2.77'
Set-TextCell $ws.Range("E40") '  -2.07%  '
Set-TextCell $ws.Range("E41") '  -0.03%  '
Set-TextCell $ws.Range("D42") '1.999.81'
Set-TextCell $ws.Range("E42") '  -0.08%  '
Set-TextCell $ws.Range("E43") '  +0.32%  '
Set-TextCell $ws.Range("D45") '10.10'
Set-TextCell $ws.Range("E45") '  -1.49%  '
Set-TextCell $ws.Range("D46") '17.47'
Set-TextCell $ws.Range("E46") '  -0.56%  '
Set-TextCell $ws.Range("E47") '  -0.78%  '
Set-TextCell $ws.Range("D48") '54.94'
Set-TextCell $ws.Range("E48") '  -1.39%  '
Set-TextCell $ws.Range("D49") '74.86'
Set-TextCell $ws.Range("E49") '  +6.50%  '
Set-TextCell $ws.Range("D50") '2.540.19'
Set-TextCell $ws.Range("E50") '  +0.89%  '
Set-TextCell $ws.Range("E51") '  +1.35%  '
